# Adds a new "2022-Q3" sheet (with fund-holding detail data) right after
# the "总计" summary sheet and before the existing "2022-Q2" sheet, and
# updates the "总计" sheet with a new top data row summarising 2022-Q3,
# shifting all the previously-existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell as TEXT (never auto-coerced to a
# number / date by Excel) while leaving the cell's style index exactly
# as it was before (no lingering "@" number-format style). We do this by
# temporarily applying a text number format, assigning the value, then
# pasting-special just the *formats* back in from a pristine, never
# touched cell so the visible/storage style reverts to the default.
# ---------------------------------------------------------------------
function Set-TextValue($ws, $addr, $val, $blankAddr) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($blankAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# =======================================================================
# 1) Insert the new "2022-Q3" worksheet, positioned right before "2022-Q2"
# =======================================================================
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$ws = $wb.Worksheets.Add($beforeSheet)
$ws.Name = "2022-Q3"

# Re-fetch sheet references AFTER the sheet collection changed so Copy /
# PasteSpecial operate against live (non-stale) objects.
$srcQ2 = $wb.Worksheets.Item("2022-Q2")

# ---- bring over the look & feel (cell styles) from the 2022-Q2 sheet ----
# header row formatting (bold / bordered / centered style)
$srcQ2.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

# data-row formatting (row index column + plain cells), replicated for
# each of the 3 data rows this sheet needs
$srcQ2.Range("A2:H2").Copy()
$ws.Range("A2:H2").PasteSpecial(-4122)
$ws.Range("A3:H3").PasteSpecial(-4122)
$ws.Range("A4:H4").PasteSpecial(-4122)

# ---- header values ----
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# ---- row 2 ----
$ws.Range("A2").Value = 0
Set-TextValue $ws "B2" "320022" "Z1"
$ws.Range("C2").Value = "诺安研究精选股票"
Set-TextValue $ws "D2" "6.17" "Z1"
Set-TextValue $ws "E2" "92.67" "Z1"
Set-TextValue $ws "F2" "4.04" "Z1"
Set-TextValue $ws "G2" "0.2493" "Z1"
$ws.Range("H2").Value = 2

# ---- row 3 ----
$ws.Range("A3").Value = 1
Set-TextValue $ws "B3" "001270" "Z1"
$ws.Range("C3").Value = "英大灵活配置混合A"
Set-TextValue $ws "D3" "0.28" "Z1"
Set-TextValue $ws "E3" "93.98" "Z1"
Set-TextValue $ws "F3" "1.74" "Z1"
Set-TextValue $ws "G3" "0.0049" "Z1"
$ws.Range("H3").Value = 10

# ---- row 4 ----
$ws.Range("A4").Value = 2
Set-TextValue $ws "B4" "001271" "Z1"
$ws.Range("C4").Value = "英大灵活配置混合B"
Set-TextValue $ws "D4" "0.28" "Z1"
Set-TextValue $ws "E4" "93.98" "Z1"
Set-TextValue $ws "F4" "1.74" "Z1"
Set-TextValue $ws "G4" "0.0049" "Z1"
$ws.Range("H4").Value = 10

# =======================================================================
# 2) Update the "总计" summary sheet: add a new 2022-Q3 row at the top of
#    the data (row 2) and push the previously existing rows down by one.
# =======================================================================
$wsTotal = $wb.Worksheets.Item(1)

# extend the row-2 style (index column + plain cells) down onto the new
# row 9 by copying row 8's formatting (keeps every style index identical
# to what Excel would have produced by inserting a row)
$wsTotal.Range("A8:D8").Copy()
$wsTotal.Range("A9:D9").PasteSpecial(-4122)

# shift the existing 7 data rows (2..8) down to (3..9), bottom row first
for ($r = 8; $r -ge 2; $r--) {
    $dst = $r + 1
    $wsTotal.Cells.Item($dst, 2).Value = $wsTotal.Cells.Item($r, 2).Value2
    $wsTotal.Cells.Item($dst, 3).Value = $wsTotal.Cells.Item($r, 3).Value2
    $wsTotal.Cells.Item($dst, 4).Value = $wsTotal.Cells.Item($r, 4).Value2
}

# the row-index column (A) is just the constant sequence 0..7
$wsTotal.Range("A9").Value = 7

# write the brand new 2022-Q3 summary row at the top (row 2)
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.26

Write-Host "2022-Q3 sheet added and 总计 sheet updated"
